# Avancement.xlsx update — "Research done, it juste to be make clean notebook"
# Updates progress percentages and comments on the Dashboard sheet to
# reflect finished NLP / computer-vision research tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Collecte de données ---------------------------------------------------
$ws.Range("D5").Value = "ETL fonctionnel"

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "Tri plus approfondi possible"

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "OK"

$ws.Range("D8").Value = "Archivé et daté"

# --- Pré-traitement texte / images ------------------------------------------
$ws.Range("D10").Value = "lowercase, retrait des mot parasite (stop_word)"

$ws.Range("D11").Value = "tokenisation et lemmatization"

$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "cv2 blur"

$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "cv2 equalizeHist"

$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "TFIDvectorizer /BOW stop_word lemma spacy ET sift + kmeans (bag of visual word)"

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "SIFT"

$ws.Range("C16").Value = 0.8
# D16 gets rich text: a plain prefix followed by a red "FAIRE DES PIPELINE" run.
$d16Prefix = "Nettoyage text done NLP création de pipeline et CV en cours "
$d16Suffix = "FAIRE DES PIPELINE"
$ws.Range("D16").Value = $d16Prefix + $d16Suffix
$d16Run = $ws.Range("D16").Characters($d16Prefix.Length + 1, $d16Suffix.Length)
$d16Run.Font.Color = 255
$d16Run.Font.Name = "Montserrat"
$d16Run.Font.Size = 11

# --- Réduction de dimension --------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "NLP :NMF voir pour LSA/PLSA/LDA - CV : PCA sur SIFT  "

$ws.Rows(19).RowHeight = 99
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Value = "NLP : Reduction du temps de calcul et surtout création de la matrice sujet-terme par factorisation"

$ws.Range("C20").Value = 0.5
$ws.Range("D20").Value = "NLP : T-SNE à 2D (visualisation) CV : T-SNE également"

$ws.Range("D21").Value = "Partie a faire avec Flask / AmChart"

# --- Visualisation -----------------------------------------------------------
$ws.Range("C22").Value = 1

$ws.Range("C23").Value = 0.5
$ws.Range("D23").Value = "utilisation après coup de AmChart ou si pas de temps Voilà"

$ws.Range("C24").Value = 0.5

$ws.Range("C25").Value = 0.5
